$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2581.818
$ws.Range("I51").Value = 1750
$ws.Range("J51").Value = 2766.6667
$ws.Range("K51").Value = 1750
$ws.Range("L51").Value = 2766.6667
$ws.Range("M51").Value = -1266
$ws.Range("N51").Value = -3734.6667

$ws.Range("H98").Value = 1009.61536
$ws.Range("I98").Value = 864.75
$ws.Range("J98").Value = 2748
$ws.Range("K98").Value = 864.75
$ws.Range("L98").Value = 2748
$ws.Range("M98").Value = 633.25
$ws.Range("N98").Value = -5744

$ws.Range("H106").Value = 9095276
$ws.Range("I106").Value = 11115726
$ws.Range("J106").Value = 3250
$ws.Range("K106").Value = 11115726
$ws.Range("L106").Value = 3250
$ws.Range("M106").Value = -11115095
$ws.Range("N106").Value = -4512

$ws.Range("H122").Value = 1009.61536
$ws.Range("I122").Value = 864.75
$ws.Range("J122").Value = 2748
$ws.Range("K122").Value = 2594.25
$ws.Range("L122").Value = 8244
$ws.Range("M122").Value = -144.25
$ws.Range("N122").Value = -13144

$ws.Range("H131").Value = 3110
$ws.Range("I131").Value = 3586.3333
$ws.Range("J131").Value = 2931.375
$ws.Range("K131").Value = 10758.9999
$ws.Range("L131").Value = 8794.125
$ws.Range("M131").Value = -5718.999899999999
$ws.Range("N131").Value = -18874.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13891656
$ws.Range("I74").Value = 1692.8
$ws.Range("J74").Value = 31254110
$ws.Range("K74").Value = 1692.8
$ws.Range("L74").Value = 31254110
$ws.Range("M74").Value = -818.8
$ws.Range("N74").Value = -31255858

$ws.Range("H77").Value = 13891656
$ws.Range("I77").Value = 1692.8
$ws.Range("J77").Value = 31254110
$ws.Range("K77").Value = 8464
$ws.Range("L77").Value = 156270550
$ws.Range("M77").Value = -4096
$ws.Range("N77").Value = -156279286

$ws.Range("H128").Value = 39958
$ws.Range("J128").Value = 39958
$ws.Range("L128").Value = 39958
$ws.Range("N128").Value = -49918

$ws.Range("H132").Value = 2408610
$ws.Range("I132").Value = 4309.2354
$ws.Range("J132").Value = 5133484
$ws.Range("K132").Value = 12927.7062
$ws.Range("L132").Value = 15400452
$ws.Range("M132").Value = -10397.7062
$ws.Range("N132").Value = -15405512

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 11065.714
$ws.Range("J57").Value = 11065.714
$ws.Range("L57").Value = 11065.714
$ws.Range("N57").Value = -12505.714

$ws.Range("H86").Value = 1869.675
$ws.Range("I86").Value = 1667.6129
$ws.Range("J86").Value = 2565.6667
$ws.Range("K86").Value = 1667.6129
$ws.Range("L86").Value = 2565.6667
$ws.Range("M86").Value = -544.6129000000001
$ws.Range("N86").Value = -4811.6667

$ws.Range("H89").Value = 1869.675
$ws.Range("I89").Value = 1667.6129
$ws.Range("J89").Value = 2565.6667
$ws.Range("K89").Value = 8338.0645
$ws.Range("L89").Value = 12828.3335
$ws.Range("M89").Value = -2722.0645
$ws.Range("N89").Value = -24060.3335

$ws.Range("H105").Value = 2720
$ws.Range("I105").Value = 2720
$ws.Range("K105").Value = 2720
$ws.Range("M105").Value = -973

$ws.Range("H134").Value = 6400
$ws.Range("I134").Value = 10000
$ws.Range("K134").Value = 30000
$ws.Range("M134").Value = -27465

$ws.Range("H136").Value = 11065.714
$ws.Range("J136").Value = 11065.714
$ws.Range("L136").Value = 11065.714
$ws.Range("N136").Value = -21265.714

$ws.Range("H138").Value = 42588
$ws.Range("I138").Value = 10000
$ws.Range("J138").Value = 50735
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 50735
$ws.Range("M138").Value = -4860
$ws.Range("N138").Value = -61015

$ws.Range("H139").Value = 69794.2
$ws.Range("J139").Value = 69794.2
$ws.Range("L139").Value = 69794.2
$ws.Range("N139").Value = -80074.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 22287.428
$ws.Range("J4").Value = 22287.428
$ws.Range("L4").Value = 22287.428
$ws.Range("N4").Value = -22511.428

$ws.Range("H58").Value = 3617.652
$ws.Range("I58").Value = 3847.7646
$ws.Range("J58").Value = 2965.6667
$ws.Range("K58").Value = 3847.7646
$ws.Range("L58").Value = 2965.6667
$ws.Range("M58").Value = -3644.7646
$ws.Range("N58").Value = -3371.6667

$ws.Range("H99").Value = 2021.3414
$ws.Range("I99").Value = 1601
$ws.Range("K99").Value = 1601
$ws.Range("M99").Value = -103

$ws.Range("H126").Value = 2021.3414
$ws.Range("I126").Value = 1601
$ws.Range("K126").Value = 4803
$ws.Range("M126").Value = -2333

$ws.Range("H136").Value = 3617.652
$ws.Range("I136").Value = 3847.7646
$ws.Range("J136").Value = 2965.6667
$ws.Range("K136").Value = 11543.2938
$ws.Range("L136").Value = 8897.000100000001
$ws.Range("M136").Value = -8993.2938
$ws.Range("N136").Value = -13997.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4244489
$ws.Range("J4").Value = 2620.9167
$ws.Range("L4").Value = 7862.750100000001
$ws.Range("N4").Value = -8086.750100000001

$ws.Range("H125").Value = 3207.5715
$ws.Range("J125").Value = 3207.5715
$ws.Range("L125").Value = 9622.7145
$ws.Range("N125").Value = -19462.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8001
$ws.Range("J2").Value = 8001
$ws.Range("L2").Value = 8001
$ws.Range("N2").Value = -8225

$ws.Range("H7").Value = 6509.9443
$ws.Range("I7").Value = 6567.615
$ws.Range("K7").Value = 6567.615
$ws.Range("M7").Value = -6455.615

$ws.Range("H94").Value = 33000
$ws.Range("J94").Value = 33000
$ws.Range("L94").Value = 33000
$ws.Range("N94").Value = -34352

$ws.Range("H126").Value = 6509.9443
$ws.Range("I126").Value = 6567.615
$ws.Range("K126").Value = 19702.845
$ws.Range("M126").Value = -17232.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 111608.25
$ws.Range("I62").Value = 145399.86
$ws.Range("K62").Value = 145399.86
$ws.Range("M62").Value = -144775.86

$ws.Range("H65").Value = 111608.25
$ws.Range("I65").Value = 145399.86
$ws.Range("K65").Value = 726999.2999999999
$ws.Range("M65").Value = -723879.2999999999

$ws.Range("H122").Value = 1281.5938
$ws.Range("I122").Value = 1247.3572
$ws.Range("K122").Value = 3742.0716
$ws.Range("M122").Value = -1292.0716

Write-Output "Applied 163 cell updates across 7 sheets"